$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 411 (shifts existing rows 411-426 down to 412-427)
$ws.Rows.Item(411).Insert()

# Populate the newly inserted row 411 with the latest weekly price record
$ws.Range("A411").Value = 4
$ws.Range("B411").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C411").Value = "Los Lagos"
$ws.Range("D411").Value = 45075
$ws.Range("E411").Value = 10
$ws.Range("F411").Value = "Fruta"
$ws.Range("G411").Value = 100108
$ws.Range("H411").Value = "Tropicales y subtropicales"
$ws.Range("I411").Value = 100108005
$ws.Range("J411").Value = "Piña"
$ws.Range("K411").Value = "Caramelo"
$ws.Range("L411").Value = "Primera"
$ws.Range("M411").Value = 80
$ws.Range("N411").Value = 18000
$ws.Range("O411").Value = 19000
$ws.Range("P411").Value = 18500
$ws.Range("Q411").Value = "`$/caja 12 unidades"
$ws.Range("R411").Value = "Ecuador"
$ws.Range("S411").Value = 1542
$ws.Range("T411").Value = 12
